# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (F) and, where a previously sold-out event is now
# on sale again, "最低票价" (G) counts/prices across all sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 8411
$ws.Cells.Item(3, 6).Value = 36780
$ws.Cells.Item(3, 7).Value = 68
$ws.Cells.Item(5, 6).Value = 623
$ws.Cells.Item(6, 6).Value = 752
$ws.Cells.Item(7, 6).Value = 472
$ws.Cells.Item(8, 6).Value = 157
$ws.Cells.Item(10, 6).Value = 826
$ws.Cells.Item(11, 6).Value = 82
$ws.Cells.Item(12, 6).Value = 668
$ws.Cells.Item(13, 6).Value = 506
$ws.Cells.Item(14, 6).Value = 31
$ws.Cells.Item(15, 6).Value = 614
$ws.Cells.Item(16, 6).Value = 174
$ws.Cells.Item(17, 6).Value = 453
$ws.Cells.Item(18, 6).Value = 437
$ws.Cells.Item(19, 6).Value = 1141
$ws.Cells.Item(21, 6).Value = 786
$ws.Cells.Item(22, 6).Value = 2455
$ws.Cells.Item(23, 6).Value = 954
$ws.Cells.Item(24, 6).Value = 537
$ws.Cells.Item(25, 6).Value = 92
$ws.Cells.Item(26, 6).Value = 1135
$ws.Cells.Item(28, 6).Value = 726
$ws.Cells.Item(29, 6).Value = 726
$ws.Cells.Item(30, 6).Value = 32
$ws.Cells.Item(31, 6).Value = 1129

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 367
$ws.Cells.Item(5, 6).Value = 325
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(11, 6).Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 599

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 599
$ws.Cells.Item(3, 6).Value = 8411
$ws.Cells.Item(5, 6).Value = 36780
$ws.Cells.Item(5, 7).Value = 68
$ws.Cells.Item(7, 6).Value = 623
$ws.Cells.Item(8, 6).Value = 752
$ws.Cells.Item(9, 6).Value = 472
$ws.Cells.Item(11, 6).Value = 157
$ws.Cells.Item(13, 6).Value = 367
$ws.Cells.Item(14, 6).Value = 325
$ws.Cells.Item(16, 6).Value = 826
$ws.Cells.Item(17, 6).Value = 82
$ws.Cells.Item(18, 6).Value = 668
$ws.Cells.Item(19, 6).Value = 506
$ws.Cells.Item(21, 6).Value = 31
$ws.Cells.Item(24, 6).Value = 5
$ws.Cells.Item(25, 6).Value = 5
$ws.Cells.Item(26, 6).Value = 614
$ws.Cells.Item(27, 6).Value = 174
$ws.Cells.Item(28, 6).Value = 453
$ws.Cells.Item(29, 6).Value = 437
$ws.Cells.Item(30, 6).Value = 1141
$ws.Cells.Item(32, 6).Value = 786
$ws.Cells.Item(33, 6).Value = 2455
$ws.Cells.Item(34, 6).Value = 954
$ws.Cells.Item(35, 6).Value = 537
$ws.Cells.Item(36, 6).Value = 92
$ws.Cells.Item(37, 6).Value = 1135
$ws.Cells.Item(40, 6).Value = 726
$ws.Cells.Item(41, 6).Value = 726
$ws.Cells.Item(42, 6).Value = 32
$ws.Cells.Item(43, 6).Value = 1129
